$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (this also updates the "Reactivos" defined name's sheet reference)
$ws.Name = "Catálogos"

# Update header/label and templating text for the remaining rows
$ws.Range("B3").Value = "{{Catalogo.Clave}}"

$ws.Range("A5").Value = "Largo"
$ws.Range("B5").Value = "{{Catalogo.Largo}}"

$ws.Range("A7").Value = "Ancho"
$ws.Range("B7").Value = "{{Catalogo.Ancho}}"

$ws.Range("A9").Value = "Activo"
$ws.Range("B9").Value = "{{Catalogo.Activo}}"

# Remove the old trailing "Activo" row (previously row 11) entirely
$ws.Rows(11).Delete()
